# Apply updated cryptocurrency price/volume figures (and one ranking swap)
# to match the data refresh from the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (leading apostrophe forces text
# interpretation even for numeric-looking strings like "1.00" or "0.999"),
# then restore the cell style to "Normal" so no stray number-format/style
# is left behind -- matches the original inline-string cells exactly.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2: Bitcoin
Set-TextValue "D2" "91.026.13"
Set-TextValue "E2" "  +3.75%  "

# Row 3: Ethereum
Set-TextValue "D3" "3.171.99"
Set-TextValue "E3" "  -0.11%  "

# Row 4: TetherUSD
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.07%  "

# Row 5: Solana
Set-TextValue "D5" "214.78"
Set-TextValue "E5" "  +3.20%  "

# Row 6: BNB
Set-TextValue "D6" "629.36"
Set-TextValue "E6" "  +3.06%  "

# Row 7: Dogecoin
Set-TextValue "D7" "0.394"
Set-TextValue "E7" "  +1.82%  "

# Row 8: XRP
Set-TextValue "D8" "0.715"
Set-TextValue "E8" "  +6.31%  "

# Row 9: USDC
Set-TextValue "D9" "1.00"
Set-TextValue "E9" "  +0.07%  "

# Row 10: LidoStakedEther
Set-TextValue "D10" "3.171.57"
Set-TextValue "E10" "  +0.05%  "

# Row 11: Cardano
Set-TextValue "D11" "0.563"
Set-TextValue "E11" "  +4.50%  "

# Row 12: TRON
Set-TextValue "E12" "  +2.17%  "

# Row 13: ShibaInu
Set-TextValue "D13" "0.0000252"
Set-TextValue "E13" "  +3.33%  "

# Row 14: WrappedBTC
Set-TextValue "D14" "90.606.56"
Set-TextValue "E14" "  +3.37%  "

# Row 15: Toncoin
Set-TextValue "D15" "5.30"
Set-TextValue "E15" "  +0.62%  "

# Row 16: WrappedliquidstakedEther2.0
Set-TextValue "D16" "3.765.14"
Set-TextValue "E16" "  +0.13%  "

# Row 17: Avalanche
Set-TextValue "D17" "32.45"
Set-TextValue "E17" "  +0.66%  "

# Row 18: WrappedEther
Set-TextValue "E18" "  -0.26%  "

# Row 19: SuiNetwork
Set-TextValue "D19" "3.29"
Set-TextValue "E19" "  +2.87%  "

# Row 20: PEPE
Set-TextValue "D20" "0.0000210"
Set-TextValue "E20" "  +58.29%  "

# Row 21: BitcoinCash
Set-TextValue "D21" "437.59"
Set-TextValue "E21" "  +5.88%  "

# Row 22: Chainlink
Set-TextValue "D22" "13.24"
Set-TextValue "E22" "  -1.92%  "

# Row 23: Uniswap
Set-TextValue "D23" "8.44"
Set-TextValue "E23" "  -0.52%  "

# Row 24: Polkadot
Set-TextValue "D24" "4.97"
Set-TextValue "E24" "  -1.96%  "

# Row 25: NEARProtocol
Set-TextValue "D25" "5.18"
Set-TextValue "E25" "  -1.48%  "

# Row 26: Aptos
Set-TextValue "D26" "11.61"
Set-TextValue "E26" "  -4.96%  "

# Row 27: Litecoin
Set-TextValue "D27" "80.28"
Set-TextValue "E27" "  +9.06%  "

# Row 28: WrappedeETH
Set-TextValue "D28" "3.342.01"
Set-TextValue "E28" "  -0.06%  "

# Row 29: Dai
Set-TextValue "E29" "  +0.08%  "

# Row 30: Binance-PegBSC-USD
Set-TextValue "D30" "0.999"
Set-TextValue "E30" "  +0.13%  "

# Row 31: Cronos
Set-TextValue "D31" "0.155"
Set-TextValue "E31" "  -4.14%  "

# Row 32: dogwifhat
Set-TextValue "D32" "3.98"
Set-TextValue "E32" "  +30.83%  "

# Row 33: InternetComputer(DFINITY)
Set-TextValue "D33" "8.27"
Set-TextValue "E33" "  +0.47%  "

# Row 34: Bittensor
Set-TextValue "D34" "518.16"
Set-TextValue "E34" "  -5.03%  "

# Row 35: RenderToken
Set-TextValue "D35" "6.91"
Set-TextValue "E35" "  -0.09%  "

# Row 36: PancakeSwap
Set-TextValue "E36" "  +0.99%  "

# Row 37: Fetch.AI
Set-TextValue "D37" "1.28"
Set-TextValue "E37" "  -3.50%  "

# Row 38: EthereumClassic
Set-TextValue "D38" "22.27"
Set-TextValue "E38" "  +1.77%  "

# Row 39: WhiteBITCoin
Set-TextValue "E39" "  +2.44%  "

# Row 40: FirstDigitalUSD
Set-TextValue "D40" "1.00"
Set-TextValue "E40" "  +0.29%  "

# Row 41: Kaspa
Set-TextValue "D41" "0.125"
Set-TextValue "E41" "  -3.58%  "

# Row 42: USDe
Set-TextValue "E42" "  -0.05%  "

# Row 43: Stacks
Set-TextValue "E43" "  -0.10%  "

# Row 44: PolygonEcosystemToken
Set-TextValue "E44" "  -1.62%  "

# Row 45: Monero
Set-TextValue "D45" "146.68"
Set-TextValue "E45" "  -1.16%  "

# Row 46: OKB
Set-TextValue "D46" "43.81"
Set-TextValue "E46" "  +1.32%  "

# Row 47: Aave
Set-TextValue "D47" "170.67"
Set-TextValue "E47" "  -1.92%  "

# Row 48: Stellar
Set-TextValue "E48" "  +0.36%  "

# Row 49: Mantle
Set-TextValue "E49" "  +7.05%  "

# Row 50: InjectiveProtocol
Set-TextValue "D50" "24.48"
Set-TextValue "E50" "  +2.46%  "

# Row 51: ARBITRUM -> ImmutableX
Set-TextValue "B51" "ImmutableX"
Set-TextValue "C51" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D51" "1.21"
Set-TextValue "E51" "  -1.98%  "
